$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last existing data row is 269 (date serial 44343 = 2021-05-27).
# New rows 270..301 continue the series up to serial 44375 (2021-06-28),
# all with 0 in columns B, C, D, matching the formatting of row 269.

$lastRow = 269
$lastSerial = 44343
$newLastRow = 301

$srcRow = $ws.Range("A" + $lastRow + ":D" + $lastRow)

for ($r = $lastRow + 1; $r -le $newLastRow; $r++) {
    $serial = $lastSerial + ($r - $lastRow)

    $destRow = $ws.Range("A" + $r + ":D" + $r)
    $srcRow.Copy($destRow)

    $ws.Range("A" + $r).Value = $serial
    $ws.Range("B" + $r).Value = 0
    $ws.Range("C" + $r).Value = 0
    $ws.Range("D" + $r).Value = 0
}
